# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.411.03"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +3.55%  "
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.839.09"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +3.61%  "
$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.031"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  +2.97%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "317.72"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("E6").Value = "  +2.43%  "
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4355"
$cell.Style = $origStyle
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3721"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +2.28%  "
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07344"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +2.64%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.8749"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +3.22%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.35"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +4.17%  "
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.985.79"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +12.00%  "
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.475"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +4.08%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.680"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +3.51%  "
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07153"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +4.16%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "82.18"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +4.07%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.032"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +2.83%  "
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.000008989"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +3.93%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.024"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +2.48%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.39"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +2.85%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.433.62"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +3.59%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.249"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +2.58%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.14"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +0.40%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.177.86"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +9.38%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "156.92"
$cell.Style = $origStyle
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.903"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +1.67%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.54"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +2.84%  "
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.259"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +3.48%  "
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.917"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +5.91%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "115.43"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +1.48%  "
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09019"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +1.04%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.199"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  +6.16%  "
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7582"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  +3.91%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.472"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  +3.32%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.866"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +4.44%  "
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.029"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +2.96%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.152"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +4.31%  "
$ws.Range("E38").Value = "  +3.42%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.05241"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.798"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  +6.31%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5143"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +4.23%  "
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1660"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  +2.86%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.521"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +3.09%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.455"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  +5.27%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "108.05"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("E46").Value = "  +4.02%  "
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.028"
$cell.Style = $origStyle
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.671"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4625"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +3.21%  "
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.884"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +8.84%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06292"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +1.35%  "
